$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

# --- Fix two player name typos/updates in existing rows ---
$ws.Range("B187").Value = "Lance McCullers Jr."
$ws.Range("B189").Value = "Peter Alonso"

# --- Append newly drafted players (rows 213-252) ---
$ws.Range("A213").Value = "bears"
$ws.Range("B213").Value = "Manny Machado"
$ws.Range("C213").Value = 29
$ws.Range("D213").Value = "SS"
$ws.Range("E213").Value = 43836

$ws.Range("A214").Value = "bellevegas"
$ws.Range("B214").Value = "Alex Reyes"
$ws.Range("C214").Value = 4
$ws.Range("D214").Value = "P"
$ws.Range("E214").Value = 43836

$ws.Range("A215").Value = "isotopes"
$ws.Range("B215").Value = "Freddie Freeman"
$ws.Range("C215").Value = 34
$ws.Range("D215").Value = "1B"
$ws.Range("E215").Value = 43836

$ws.Range("A216").Value = "ds9"
$ws.Range("B216").Value = "Max Scherzer"
$ws.Range("C216").Value = 39
$ws.Range("D216").Value = "P"
$ws.Range("E216").Value = 43836

$ws.Range("A217").Value = "dembums"
$ws.Range("B217").Value = "Aaron Nola"
$ws.Range("C217").Value = 30
$ws.Range("D217").Value = "P"
$ws.Range("E217").Value = 43836

$ws.Range("A218").Value = "pkdodgers"
$ws.Range("B218").Value = "Jesse Winker"
$ws.Range("C218").Value = 6
$ws.Range("D218").Value = "OF"
$ws.Range("E218").Value = 43836

$ws.Range("A219").Value = "bellevegas"
$ws.Range("B219").Value = "Michael Brantley"
$ws.Range("C219").Value = 19
$ws.Range("D219").Value = "OF"
$ws.Range("E219").Value = 43836

$ws.Range("A220").Value = "pasadena"
$ws.Range("B220").Value = "Jose Altuve"
$ws.Range("C220").Value = 25
$ws.Range("D220").Value = "2B"
$ws.Range("E220").Value = 43836

$ws.Range("A221").Value = "deener"
$ws.Range("B221").Value = "Yasmani Grandal"
$ws.Range("C221").Value = 23
$ws.Range("D221").Value = "C"
$ws.Range("E221").Value = 43837

$ws.Range("A222").Value = "drjames"
$ws.Range("B222").Value = "Kris Bryant"
$ws.Range("C222").Value = 29
$ws.Range("D222").Value = "OF"
$ws.Range("E222").Value = 43837

$ws.Range("A223").Value = "rippe"
$ws.Range("B223").Value = "Garrett Hampson"
$ws.Range("C223").Value = 10
$ws.Range("D223").Value = "2B"
$ws.Range("E223").Value = 43837

$ws.Range("A224").Value = "rippe"
$ws.Range("B224").Value = "Willson Contreras"
$ws.Range("C224").Value = 18
$ws.Range("D224").Value = "C"
$ws.Range("E224").Value = 43837

$ws.Range("A225").Value = "bellevegas"
$ws.Range("B225").Value = "Zack Greinke"
$ws.Range("C225").Value = 23
$ws.Range("D225").Value = "P"
$ws.Range("E225").Value = 43837

$ws.Range("A226").Value = "bellevegas"
$ws.Range("B226").Value = "Clayton Kershaw"
$ws.Range("C226").Value = 28
$ws.Range("D226").Value = "P"
$ws.Range("E226").Value = 43837

$ws.Range("A227").Value = "deano"
$ws.Range("B227").Value = "Buster Posey"
$ws.Range("C227").Value = 8
$ws.Range("D227").Value = "C"
$ws.Range("E227").Value = 43837

$ws.Range("A228").Value = "dsb"
$ws.Range("B228").Value = "Chris Sale"
$ws.Range("C228").Value = 33
$ws.Range("D228").Value = "P"
$ws.Range("E228").Value = 43837

$ws.Range("A229").Value = "deener"
$ws.Range("B229").Value = "Zach Davies"
$ws.Range("C229").Value = 1
$ws.Range("D229").Value = "P"
$ws.Range("E229").Value = 43838

$ws.Range("A230").Value = "chicago"
$ws.Range("B230").Value = "Starling Marte"
$ws.Range("C230").Value = 31
$ws.Range("D230").Value = "OF"
$ws.Range("E230").Value = 43838

$ws.Range("A231").Value = "bears"
$ws.Range("B231").Value = "Paul Goldschmidt"
$ws.Range("C231").Value = 24
$ws.Range("D231").Value = "1B"
$ws.Range("E231").Value = 43838

$ws.Range("A232").Value = "ds9"
$ws.Range("B232").Value = "Bryce Harper"
$ws.Range("C232").Value = 34
$ws.Range("D232").Value = "OF"
$ws.Range("E232").Value = 43838

$ws.Range("A233").Value = "chicago"
$ws.Range("B233").Value = "Nomar Mazara"
$ws.Range("C233").Value = 16
$ws.Range("D233").Value = "OF"
$ws.Range("E233").Value = 43838

$ws.Range("A234").Value = "bellevegas"
$ws.Range("B234").Value = "Starlin Castro"
$ws.Range("C234").Value = 8
$ws.Range("D234").Value = "MI"
$ws.Range("E234").Value = 43838

$ws.Range("A235").Value = "sturgeon"
$ws.Range("B235").Value = "Yu Darvish"
$ws.Range("C235").Value = 24
$ws.Range("D235").Value = "P"
$ws.Range("E235").Value = 43838

$ws.Range("A236").Value = "bellevegas"
$ws.Range("B236").Value = "Dakota Hudson"
$ws.Range("C236").Value = 8
$ws.Range("D236").Value = "P"
$ws.Range("E236").Value = 43838

$ws.Range("A237").Value = "marmaduke"
$ws.Range("B237").Value = "Giancarlo Stanton"
$ws.Range("C237").Value = 34
$ws.Range("D237").Value = "OF"
$ws.Range("E237").Value = 43839

$ws.Range("A238").Value = "deano"
$ws.Range("B238").Value = "Anthony Rizzo"
$ws.Range("C238").Value = 29
$ws.Range("D238").Value = "1B"
$ws.Range("E238").Value = 43839

$ws.Range("A239").Value = "marmaduke"
$ws.Range("B239").Value = "Yadier Molina"
$ws.Range("C239").Value = 8
$ws.Range("D239").Value = "C"
$ws.Range("E239").Value = 43839

$ws.Range("A240").Value = "pkdodgers"
$ws.Range("B240").Value = "Kenley Jansen"
$ws.Range("C240").Value = 20
$ws.Range("D240").Value = "P"
$ws.Range("E240").Value = 43839

$ws.Range("A241").Value = "dsb"
$ws.Range("B241").Value = "Josh Donaldson"
$ws.Range("C241").Value = 24
$ws.Range("D241").Value = "CI"
$ws.Range("E241").Value = 43839

$ws.Range("A242").Value = "deano"
$ws.Range("B242").Value = "Wilson Ramos"
$ws.Range("C242").Value = 13
$ws.Range("D242").Value = "C"
$ws.Range("E242").Value = 43839

$ws.Range("A243").Value = "sturgeon"
$ws.Range("B243").Value = "Aroldis Chapman"
$ws.Range("C243").Value = 26
$ws.Range("D243").Value = "P"
$ws.Range("E243").Value = 43839

$ws.Range("A244").Value = "balco"
$ws.Range("B244").Value = "Zack Wheeler"
$ws.Range("C244").Value = 20
$ws.Range("D244").Value = "P"
$ws.Range("E244").Value = 43839

$ws.Range("A245").Value = "pkdodgers"
$ws.Range("B245").Value = "Jorge Polanco"
$ws.Range("C245").Value = 15
$ws.Range("D245").Value = "SS"
$ws.Range("E245").Value = 43840

$ws.Range("A246").Value = "pasadena"
$ws.Range("B246").Value = "Andrew Heaney"
$ws.Range("C246").Value = 16
$ws.Range("D246").Value = "P"
$ws.Range("E246").Value = 43840

$ws.Range("A247").Value = "chicago"
$ws.Range("B247").Value = "Jason Heyward"
$ws.Range("C247").Value = 7
$ws.Range("D247").Value = "OF"
$ws.Range("E247").Value = 43840

$ws.Range("A248").Value = "dsb"
$ws.Range("B248").Value = "Corey Kluber"
$ws.Range("C248").Value = 21
$ws.Range("D248").Value = "P"
$ws.Range("E248").Value = 43840

$ws.Range("A249").Value = "sturgeon"
$ws.Range("B249").Value = "Raisel Iglesias"
$ws.Range("C249").Value = 22
$ws.Range("D249").Value = "P"
$ws.Range("E249").Value = 43840

$ws.Range("A250").Value = "marmaduke"
$ws.Range("B250").Value = "Travis Shaw"
$ws.Range("C250").Value = 6
$ws.Range("D250").Value = "3B"
$ws.Range("E250").Value = 43840

$ws.Range("A251").Value = "bears"
$ws.Range("B251").Value = "Luis Arraez"
$ws.Range("C251").Value = 9
$ws.Range("D251").Value = "MI"
$ws.Range("E251").Value = 43840

$ws.Range("A252").Value = "marmaduke"
$ws.Range("B252").Value = "Carlos Santana"
$ws.Range("C252").Value = 19
$ws.Range("D252").Value = "1B"
$ws.Range("E252").Value = 43840

# --- Restore view/selection state ---
$ws.Range("B215").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 199
$win.ScrollColumn = 1